# The "label boolean" column setup is gone: columns must now be entered
# directly by their letters. Insert two new columns (E: "Prénom", F:
# "Note/10,00") before the existing "Temps utilisé" column on Feuille5,
# which shifts that column from E to G, and fill in the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuille5")

# Insert two blank columns at E:F — this pushes the old E ("Temps utilisé")
# column (and its formatting) to G, exactly matching the target layout.
$ws.Range("E1:F1").EntireColumn.Insert()

# Header row
$ws.Range("E1").Value = "Prénom"
$ws.Range("F1").Value = "Note/10,00"

# Data rows: new "Prénom" (E) and "Note/10,00" (F) values per student
$ws.Range("E2").Value = "Houzefa"
$ws.Range("F2").Value = "7,83"

$ws.Range("E3").Value = "Yasmine"
$ws.Range("F3").Value = "7,83"

$ws.Range("E4").Value = "Zina"
$ws.Range("F4").Value = "7,28"

$ws.Range("E5").Value = "Aboubaker"
$ws.Range("F5").Value = "7,98"

$ws.Range("E6").Value = "Yasmine"
$ws.Range("F6").Value = "7,83"

$ws.Range("E7").Value = "Hassan Mahamat"
$ws.Range("F7").Value = "7,52"

$ws.Range("E8").Value = "Yacine"
$ws.Range("F8").Value = "8,07"

$ws.Range("E9").Value = "Paola"
$ws.Range("F9").Value = "7,16"

$ws.Range("E10").Value = "Rodolphe"
$ws.Range("F10").Value = "8,88"

$ws.Range("E11").Value = "Nouh"
$ws.Range("F11").Value = "7,35"

$ws.Range("E12").Value = "Iness"
$ws.Range("F12").Value = "6,51"

$ws.Range("E13").Value = "Zakaria"
$ws.Range("F13").Value = "7,70"

$ws.Range("E14").Value = "Christian"
$ws.Range("F14").Value = "6,72"
